# "Error Calculations and Plots"
#
# This missing-data sample dropped two records entirely (RM 232 and
# SC 92), shifting every row below them up by one/two, and the
# simulated "missing" (B/C/D/E/F) cells were re-drawn for several
# rows. Below: remove the two now-empty trailing rows, then restate
# every row whose values actually differ from the original.
#
# Helper: writes a blank-but-text placeholder cell (matches the
# workbook's existing convention for "missing" values — an empty
# string, not a truly blank/removed cell) without leaving Excel's
# quote-prefix formatting behind.
function Set-BlankMarker($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RM 232 (old row 34) and SC 232 (old row 35) no longer belong in the
# table — drop them so the sheet ends at row 33.
$ws.Range("A34:A35").EntireRow.Delete()

# Row 2 (RM 2): B was re-imputed
$ws.Range("C2").Value = 14.9

# Row 3 (RM 8): B now missing
Set-BlankMarker "C3"

# Row 4 (RM 9): B now missing
Set-BlankMarker "C4"

# Row 5 (RM 14): F now missing
Set-BlankMarker "F5"

# Row 8 (RM 38): F re-imputed
$ws.Range("F8").Value = 17.05

# Row 10 (RM 52 a): F re-imputed
$ws.Range("F10").Value = 16.43

# Row 11 (RM 58): B re-imputed
$ws.Range("C11").Value = 11.4

# Row 12 (RM 81): F now missing
Set-BlankMarker "F12"

# Row 13 (RM 88): B now missing
Set-BlankMarker "C13"

# Row 15 (RM 95): F re-imputed
$ws.Range("F15").Value = 16.2

# Row 18 (RM 120): F now missing
Set-BlankMarker "F18"

# Row 19 (RM 125): F now missing
Set-BlankMarker "F19"

# Row 21 (RM 135): B re-imputed
$ws.Range("C21").Value = 12.7

# Row 25 (RM 145): B now missing, F re-imputed
Set-BlankMarker "C25"
$ws.Range("F25").Value = 16.6

# Row 26: was RM 232 -> now holds the old SC 5 record (row shifted up by one)
$ws.Range("A26").Value = "SC 5"
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5.0
$ws.Range("F26").Value = 17.38

# Row 27: was SC 5 -> now holds the old SC 101 record
$ws.Range("A27").Value = "SC 101"
$ws.Range("C27").Value = 10.0
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10.0
$ws.Range("F27").Value = 17.0

# Row 28: was SC 92 (removed) -> now holds the old SC 105 record
$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

# Row 29: was SC 101 -> now holds the old SC 119 record (still missing F)
$ws.Range("A29").Value = "SC 119"
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13.0
$ws.Range("E29").Value = -6.8

# Row 30: was SC 105 -> now holds the old SC 120 record
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

# Row 31: was SC 119 -> now holds the old SC 132 record
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

# Row 32: was SC 120 -> now holds the old SC 193 record
$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

# Row 33: was SC 132 -> now holds the old SC 232 record (F now missing)
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
Set-BlankMarker "F33"
